$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Motor")
Write-Host $ws.Name
